$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 35996.8
$ws.Range("I62").Value = 60086.117
$ws.Range("J62").Value = 4495.385
$ws.Range("K62").Value = 60086.117
$ws.Range("L62").Value = 4495.385
$ws.Range("M62").Value = -59462.117
$ws.Range("N62").Value = -5743.385
$ws.Range("H65").Value = 35996.8
$ws.Range("I65").Value = 60086.117
$ws.Range("J65").Value = 4495.385
$ws.Range("K65").Value = 300430.585
$ws.Range("L65").Value = 22476.925
$ws.Range("M65").Value = -297310.585
$ws.Range("N65").Value = -28716.925
$ws.Range("H76").Value = 45836308
$ws.Range("I76").Value = 47829060
$ws.Range("K76").Value = 47829060
$ws.Range("M76").Value = -47828745
$ws.Range("H79").Value = 45836308
$ws.Range("I79").Value = 47829060
$ws.Range("K79").Value = 47829060
$ws.Range("M79").Value = -47827968
$ws.Range("H132").Value = 1404.7291
$ws.Range("I132").Value = 1555.4
$ws.Range("J132").Value = 999.0769
$ws.Range("K132").Value = 4666.200000000001
$ws.Range("L132").Value = 2997.2307
$ws.Range("M132").Value = -2136.200000000001
$ws.Range("N132").Value = -8057.2307
$ws.Range("H135").Value = 1128.6285
$ws.Range("I135").Value = 640.6875
$ws.Range("K135").Value = 5766.1875
$ws.Range("M135").Value = -3231.1875
$ws.Range("H138").Value = 937.5599999999999
$ws.Range("I138").Value = 558.57745
$ws.Range("J138").Value = 1865.4138
$ws.Range("K138").Value = 1675.73235
$ws.Range("L138").Value = 5596.2414
$ws.Range("M138").Value = 3464.26765
$ws.Range("N138").Value = -15876.2414
$ws.Range("H141").Value = 1687.3448
$ws.Range("I141").Value = 663.4666999999999
$ws.Range("K141").Value = 1990.4001
$ws.Range("M141").Value = 3189.5999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19665.477
$ws.Range("I32").Value = 24887.764
$ws.Range("K32").Value = 24887.764
$ws.Range("M32").Value = -24600.764
$ws.Range("H61").Value = 1065.0785
$ws.Range("I61").Value = 710.95123
$ws.Range("J61").Value = 2517
$ws.Range("K61").Value = 710.95123
$ws.Range("L61").Value = 2517
$ws.Range("M61").Value = -498.95123
$ws.Range("N61").Value = -2941
$ws.Range("H74").Value = 1044.1041
$ws.Range("I74").Value = 986.4722
$ws.Range("K74").Value = 986.4722
$ws.Range("M74").Value = -112.4722
$ws.Range("H77").Value = 1044.1041
$ws.Range("I77").Value = 986.4722
$ws.Range("K77").Value = 4932.361
$ws.Range("M77").Value = -564.3609999999999
$ws.Range("H132").Value = 961.3103599999999
$ws.Range("I132").Value = 886.38464
$ws.Range("J132").Value = 1610.6666
$ws.Range("K132").Value = 2659.15392
$ws.Range("L132").Value = 4831.9998
$ws.Range("M132").Value = -129.1539199999997
$ws.Range("N132").Value = -9891.9998
$ws.Range("H136").Value = 1065.0785
$ws.Range("I136").Value = 710.95123
$ws.Range("J136").Value = 2517
$ws.Range("K136").Value = 2132.85369
$ws.Range("L136").Value = 7551
$ws.Range("M136").Value = 417.1463100000001
$ws.Range("N136").Value = -12651

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12727.379
$ws.Range("I134").Value = 1068.3472
$ws.Range("J134").Value = 68690.734
$ws.Range("K134").Value = 3205.0416
$ws.Range("L134").Value = 206072.202
$ws.Range("M134").Value = -670.0415999999996
$ws.Range("N134").Value = -211142.202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 123.882355
$ws.Range("I7").Value = 151.90909
$ws.Range("K7").Value = 151.90909
$ws.Range("M7").Value = -38.90908999999999
$ws.Range("H31").Value = 2836.9778
$ws.Range("I31").Value = 3273.625
$ws.Range("J31").Value = 2337.9524
$ws.Range("K31").Value = 3273.625
$ws.Range("L31").Value = 2337.9524
$ws.Range("M31").Value = -2978.625
$ws.Range("N31").Value = -2927.9524
$ws.Range("H34").Value = 2836.9778
$ws.Range("I34").Value = 3273.625
$ws.Range("J34").Value = 2337.9524
$ws.Range("K34").Value = 3273.625
$ws.Range("L34").Value = 2337.9524
$ws.Range("M34").Value = -3071.625
$ws.Range("N34").Value = -2741.9524
$ws.Range("H58").Value = 2620.8596
$ws.Range("I58").Value = 825.44446
$ws.Range("J58").Value = 9353.666999999999
$ws.Range("K58").Value = 825.44446
$ws.Range("L58").Value = 9353.666999999999
$ws.Range("M58").Value = -622.44446
$ws.Range("N58").Value = -9759.666999999999
$ws.Range("H94").Value = 2711.139
$ws.Range("J94").Value = 2607.3928
$ws.Range("L94").Value = 2607.3928
$ws.Range("N94").Value = -3509.3928
$ws.Range("H97").Value = 24333.334
$ws.Range("J97").Value = 24333.334
$ws.Range("L97").Value = 24333.334
$ws.Range("N97").Value = -26315.334
$ws.Range("H99").Value = 2246.3171
$ws.Range("I99").Value = 1762.963
$ws.Range("K99").Value = 1762.963
$ws.Range("M99").Value = -264.963
$ws.Range("H126").Value = 2246.3171
$ws.Range("I126").Value = 1762.963
$ws.Range("K126").Value = 5288.889
$ws.Range("M126").Value = -2818.889
$ws.Range("H132").Value = 1067.93
$ws.Range("I132").Value = 810.24286
$ws.Range("J132").Value = 1669.2
$ws.Range("K132").Value = 2430.72858
$ws.Range("L132").Value = 5007.6
$ws.Range("M132").Value = 99.27142000000003
$ws.Range("N132").Value = -10067.6
$ws.Range("H134").Value = 976.7789299999999
$ws.Range("I134").Value = 923.2727
$ws.Range("J134").Value = 1205.6666
$ws.Range("K134").Value = 2769.8181
$ws.Range("L134").Value = 3616.9998
$ws.Range("M134").Value = -234.8181
$ws.Range("N134").Value = -8686.9998
$ws.Range("H136").Value = 2620.8596
$ws.Range("I136").Value = 825.44446
$ws.Range("J136").Value = 9353.666999999999
$ws.Range("K136").Value = 2476.33338
$ws.Range("L136").Value = 28061.001
$ws.Range("M136").Value = 73.66661999999997
$ws.Range("N136").Value = -33161.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 839.5357
$ws.Range("I5").Value = 844
$ws.Range("J5").Value = 819
$ws.Range("K5").Value = 2532
$ws.Range("L5").Value = 2457
$ws.Range("M5").Value = -2420
$ws.Range("N5").Value = -2681
$ws.Range("H56").Value = 2984.8572
$ws.Range("I56").Value = 2984.8572
$ws.Range("K56").Value = 2984.8572
$ws.Range("M56").Value = -2454.8572
$ws.Range("H68").Value = 776.7692
$ws.Range("I68").Value = 673.5
$ws.Range("J68").Value = 822.6667
$ws.Range("K68").Value = 2020.5
$ws.Range("L68").Value = 2468.0001
$ws.Range("M68").Value = -1209.5
$ws.Range("N68").Value = -4090.0001
$ws.Range("H71").Value = 776.7692
$ws.Range("I71").Value = 673.5
$ws.Range("J71").Value = 822.6667
$ws.Range("K71").Value = 6061.5
$ws.Range("L71").Value = 7404.0003
$ws.Range("M71").Value = -2005.5
$ws.Range("N71").Value = -15516.0003
$ws.Range("H107").Value = 1555843
$ws.Range("I107").Value = 216
$ws.Range("J107").Value = 3889283.5
$ws.Range("K107").Value = 648
$ws.Range("L107").Value = 11667850.5
$ws.Range("M107").Value = 1272
$ws.Range("N107").Value = -11671690.5
$ws.Range("H131").Value = 5694611
$ws.Range("J131").Value = 1472.987
$ws.Range("L131").Value = 4418.961
$ws.Range("N131").Value = -14498.961
$ws.Range("H132").Value = 737.875
$ws.Range("I132").Value = 633
$ws.Range("J132").Value = 1052.5
$ws.Range("K132").Value = 5697
$ws.Range("L132").Value = 9472.5
$ws.Range("M132").Value = -3167
$ws.Range("N132").Value = -14532.5
$ws.Range("H135").Value = 839.5357
$ws.Range("I135").Value = 844
$ws.Range("J135").Value = 819
$ws.Range("K135").Value = 7596
$ws.Range("L135").Value = 7371
$ws.Range("M135").Value = -5061
$ws.Range("N135").Value = -12441
$ws.Range("H137").Value = 28325.104
$ws.Range("I137").Value = 1627.05
$ws.Range("J137").Value = 56428.316
$ws.Range("K137").Value = 4881.15
$ws.Range("L137").Value = 169284.948
$ws.Range("M137").Value = 218.8500000000004
$ws.Range("N137").Value = -179484.948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 5008.2856
$ws.Range("I53").Value = 4635.6
$ws.Range("J53").Value = 5940
$ws.Range("K53").Value = 4635.6
$ws.Range("L53").Value = 5940
$ws.Range("M53").Value = -4004.6
$ws.Range("N53").Value = -7202
$ws.Range("H70").Value = 4350
$ws.Range("I70").Value = 4107.143
$ws.Range("K70").Value = 4107.143
$ws.Range("M70").Value = -3837.143
$ws.Range("H73").Value = 4350
$ws.Range("I73").Value = 4107.143
$ws.Range("K73").Value = 4107.143
$ws.Range("M73").Value = -3171.143
$ws.Range("H122").Value = 16229294
$ws.Range("I122").Value = 15966091
$ws.Range("J122").Value = 16667967
$ws.Range("K122").Value = 47898273
$ws.Range("L122").Value = 50003901
$ws.Range("M122").Value = -47895823
$ws.Range("N122").Value = -50008801

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1907.5834
$ws.Range("I46").Value = 831.8333
$ws.Range("J46").Value = 2983.3333
$ws.Range("K46").Value = 831.8333
$ws.Range("L46").Value = 2983.3333
$ws.Range("M46").Value = -643.8333
$ws.Range("N46").Value = -3359.3333
$ws.Range("H132").Value = 1987.2745
$ws.Range("I132").Value = 1885.4445
$ws.Range("J132").Value = 2751
$ws.Range("K132").Value = 5656.333500000001
$ws.Range("L132").Value = 8253
$ws.Range("M132").Value = -3126.333500000001
$ws.Range("N132").Value = -13313
$ws.Range("H136").Value = 1666.8
$ws.Range("I136").Value = 974.5
$ws.Range("J136").Value = 13666.667
$ws.Range("K136").Value = 2923.5
$ws.Range("L136").Value = 41000.001
$ws.Range("M136").Value = -373.5
$ws.Range("N136").Value = -46100.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1207.4
$ws.Range("I81").Value = 550.5833
$ws.Range("J81").Value = 3834.6667
$ws.Range("K81").Value = 1101.1666
$ws.Range("L81").Value = 7669.3334
$ws.Range("M81").Value = -40.16660000000002
$ws.Range("N81").Value = -9791.3334
$ws.Range("H84").Value = 1207.4
$ws.Range("I84").Value = 550.5833
$ws.Range("J84").Value = 3834.6667
$ws.Range("K84").Value = 5505.833000000001
$ws.Range("L84").Value = 38346.667
$ws.Range("M84").Value = -201.8330000000005
$ws.Range("N84").Value = -48954.667
$ws.Range("H127").Value = 57000
$ws.Range("J127").Value = 57000
$ws.Range("L127").Value = 57000
$ws.Range("N127").Value = -66920
$ws.Range("H132").Value = 548.7286
$ws.Range("I132").Value = 535.6429000000001
$ws.Range("J132").Value = 601.0714
$ws.Range("K132").Value = 1606.9287
$ws.Range("L132").Value = 1803.2142
$ws.Range("M132").Value = 923.0712999999998
$ws.Range("N132").Value = -6863.2142
